$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.130.00'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.42%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.390.85'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.13%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '568.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.46'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.23%  '

$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.610'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.30%  '

$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.387.05'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.28%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.15'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.57%  '

$ws.Range("E11").Value = '  -3.75%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.440'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.48%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.973.09'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.00%  '

$ws.Range("E14").Value = '  -0.17%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000187'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.31%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.76'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.65%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.183.00'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.41%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.385.74'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.88%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.33'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.87%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.91'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.27%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '373.76'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.51%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.95'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.91%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.547'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.31%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.998'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.53%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '71.58'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.37%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000117'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.20%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.87'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.63%  '

$ws.Range("E28").Value = '  -2.60%  '

$ws.Range("E29").Value = '  -0.11%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.46'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.35%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.04'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.31%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.01'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.30%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.12'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.75%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.12'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.03%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.58'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.11%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '159.56'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.85%  '

$ws.Range("E37").Value = '  -0.32%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0757'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.75%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '26.61'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.64%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.72'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.21%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.831.30'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.13%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.58'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.23%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.59'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.05%  '

$ws.Range("E44").Value = '  -2.93%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.760'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.61%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.55'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.55%  '

$ws.Range("E47").Value = '  -2.84%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '310.08'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.33%  '

$ws.Range("E49").Value = '  -0.26%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.53'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.10%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.853'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.17%  '
